$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(16, 3).Value = "73145412"
$ws.Cells.Item(16, 4).Value = "JOSE JESUS ORTIZ MARTINEZ"
$ws.Cells.Item(16, 5).Value = "2206"
$ws.Cells.Item(16, 6).Value = 40000

$ws.Cells.Item(17, 3).Value = "73145412"
$ws.Cells.Item(17, 4).Value = "JOSE JESUS ORTIZ MARTINEZ"
$ws.Cells.Item(17, 5).Value = "2207"
$ws.Cells.Item(17, 6).Value = 40000

$ws.Cells.Item(18, 3).Value = "73145412"
$ws.Cells.Item(18, 4).Value = "JOSE JESUS ORTIZ MARTINEZ"
$ws.Cells.Item(18, 5).Value = "2208"
$ws.Cells.Item(18, 6).Value = 40000

$ws.Cells.Item(19, 3).Value = "73145412"
$ws.Cells.Item(19, 4).Value = "JOSE JESUS ORTIZ MARTINEZ"
$ws.Cells.Item(19, 5).Value = "2209"
$ws.Cells.Item(19, 6).Value = 40000

$ws.Cells.Item(20, 3).Value = "9265520"
$ws.Cells.Item(20, 4).Value = "NAVIS JOSE CARO CARO"
$ws.Cells.Item(20, 5).Value = "2210"
$ws.Cells.Item(20, 6).Value = 40000

$ws.Cells.Item(21, 3).Value = "73145412"
$ws.Cells.Item(21, 4).Value = "JOSE JESUS ORTIZ MARTINEZ"
$ws.Cells.Item(21, 5).Value = "2210"
$ws.Cells.Item(21, 6).Value = 40000

$ws.Cells.Item(22, 3).Value = "9265520"
$ws.Cells.Item(22, 4).Value = "NAVIS JOSE CARO CARO"
$ws.Cells.Item(22, 5).Value = "2211"
$ws.Cells.Item(22, 6).Value = 40000

$ws.Cells.Item(23, 3).Value = "73145412"
$ws.Cells.Item(23, 4).Value = "JOSE JESUS ORTIZ MARTINEZ"
$ws.Cells.Item(23, 5).Value = "2211"
$ws.Cells.Item(23, 6).Value = 40000

$ws.Cells.Item(24, 3).Value = "9265520"
$ws.Cells.Item(24, 4).Value = "NAVIS JOSE CARO CARO"
$ws.Cells.Item(24, 5).Value = "2212"
$ws.Cells.Item(24, 6).Value = 40000

$ws.Cells.Item(25, 3).Value = "73145412"
$ws.Cells.Item(25, 4).Value = "JOSE JESUS ORTIZ MARTINEZ"
$ws.Cells.Item(25, 5).Value = "2212"
$ws.Cells.Item(25, 6).Value = 40000

$ws.Cells.Item(26, 3).Value = "9265520"
$ws.Cells.Item(26, 4).Value = "NAVIS JOSE CARO CARO"
$ws.Cells.Item(26, 5).Value = "2301"
$ws.Cells.Item(26, 6).Value = 40000

$ws.Cells.Item(27, 3).Value = "73145412"
$ws.Cells.Item(27, 4).Value = "JOSE JESUS ORTIZ MARTINEZ"
$ws.Cells.Item(27, 5).Value = "2301"
$ws.Cells.Item(27, 6).Value = 40000

$ws.Cells.Item(28, 3).Value = "9265520"
$ws.Cells.Item(28, 4).Value = "NAVIS JOSE CARO CARO"
$ws.Cells.Item(28, 5).Value = "2302"
$ws.Cells.Item(28, 6).Value = 40000

$ws.Cells.Item(29, 3).Value = "73145412"
$ws.Cells.Item(29, 4).Value = "JOSE JESUS ORTIZ MARTINEZ"
$ws.Cells.Item(29, 5).Value = "2302"
$ws.Cells.Item(29, 6).Value = 40000

$ws.Cells.Item(30, 3).Value = "9265520"
$ws.Cells.Item(30, 4).Value = "NAVIS JOSE CARO CARO"
$ws.Cells.Item(30, 5).Value = "2303"
$ws.Cells.Item(30, 6).Value = 40000

$ws.Cells.Item(31, 3).Value = "73145412"
$ws.Cells.Item(31, 4).Value = "JOSE JESUS ORTIZ MARTINEZ"
$ws.Cells.Item(31, 5).Value = "2303"
$ws.Cells.Item(31, 6).Value = 40000

$ws.Cells.Item(32, 3).Value = "9265520"
$ws.Cells.Item(32, 4).Value = "NAVIS JOSE CARO CARO"
$ws.Cells.Item(32, 5).Value = "2304"
$ws.Cells.Item(32, 6).Value = 40000

$ws.Cells.Item(33, 3).Value = "73145412"
$ws.Cells.Item(33, 4).Value = "JOSE JESUS ORTIZ MARTINEZ"
$ws.Cells.Item(33, 5).Value = "2304"
$ws.Cells.Item(33, 6).Value = 40000

$ws.Cells.Item(34, 3).Value = "9265520"
$ws.Cells.Item(34, 4).Value = "NAVIS JOSE CARO CARO"
$ws.Cells.Item(34, 5).Value = "2305"
$ws.Cells.Item(34, 6).Value = 40000

$ws.Cells.Item(35, 3).Value = "73145412"
$ws.Cells.Item(35, 4).Value = "JOSE JESUS ORTIZ MARTINEZ"
$ws.Cells.Item(35, 5).Value = "2305"
$ws.Cells.Item(35, 6).Value = 40000

$ws.Cells.Item(36, 3).Value = "9265520"
$ws.Cells.Item(36, 4).Value = "NAVIS JOSE CARO CARO"
$ws.Cells.Item(36, 5).Value = "2306"
$ws.Cells.Item(36, 6).Value = 40000

$ws.Cells.Item(37, 3).Value = "73145412"
$ws.Cells.Item(37, 4).Value = "JOSE JESUS ORTIZ MARTINEZ"
$ws.Cells.Item(37, 5).Value = "2306"
$ws.Cells.Item(37, 6).Value = 40000

$ws.Cells.Item(38, 3).Value = "9265520"
$ws.Cells.Item(38, 4).Value = "NAVIS JOSE CARO CARO"
$ws.Cells.Item(38, 5).Value = "2307"
$ws.Cells.Item(38, 6).Value = 40000

$ws.Cells.Item(39, 3).Value = "73145412"
$ws.Cells.Item(39, 4).Value = "JOSE JESUS ORTIZ MARTINEZ"
$ws.Cells.Item(39, 5).Value = "2307"
$ws.Cells.Item(39, 6).Value = 40000

$ws.Cells.Item(40, 3).Value = "9265520"
$ws.Cells.Item(40, 4).Value = "NAVIS JOSE CARO CARO"
$ws.Cells.Item(40, 5).Value = "2308"
$ws.Cells.Item(40, 6).Value = 40000

$ws.Cells.Item(41, 3).Value = "73145412"
$ws.Cells.Item(41, 4).Value = "JOSE JESUS ORTIZ MARTINEZ"
$ws.Cells.Item(41, 5).Value = "2308"
$ws.Cells.Item(41, 6).Value = 40000

$ws.Cells.Item(42, 3).Value = "9265520"
$ws.Cells.Item(42, 4).Value = "NAVIS JOSE CARO CARO"
$ws.Cells.Item(42, 5).Value = "2309"
$ws.Cells.Item(42, 6).Value = 40000

$ws.Cells.Item(43, 3).Value = "73145412"
$ws.Cells.Item(43, 4).Value = "JOSE JESUS ORTIZ MARTINEZ"
$ws.Cells.Item(43, 5).Value = "2309"
$ws.Cells.Item(43, 6).Value = 40000

$ws.Cells.Item(44, 3).Value = "9265520"
$ws.Cells.Item(44, 4).Value = "NAVIS JOSE CARO CARO"
$ws.Cells.Item(44, 5).Value = "2310"
$ws.Cells.Item(44, 6).Value = 40000

$ws.Cells.Item(45, 3).Value = "73145412"
$ws.Cells.Item(45, 4).Value = "JOSE JESUS ORTIZ MARTINEZ"
$ws.Cells.Item(45, 5).Value = "2310"
$ws.Cells.Item(45, 6).Value = 40000

$ws.Cells.Item(46, 3).Value = "9265520"
$ws.Cells.Item(46, 4).Value = "NAVIS JOSE CARO CARO"
$ws.Cells.Item(46, 5).Value = "2311"
$ws.Cells.Item(46, 6).Value = 37333

$ws.Cells.Item(47, 3).Value = "73145412"
$ws.Cells.Item(47, 4).Value = "JOSE JESUS ORTIZ MARTINEZ"
$ws.Cells.Item(47, 5).Value = "2311"
$ws.Cells.Item(47, 6).Value = 37333
